$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp in A1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 24 de Junio de 2020 a las 13:48"

# Refreshed country case-count data (re-sorted by "Casos totales" desc, so some
# neighbouring countries swap rows where their totals crossed over)
$ws.Cells.Item(5,2).Value = 1152066
$ws.Cells.Item(5,3).Value = 587
$ws.Cells.Item(5,5).Value = 485933
$ws.Cells.Item(5,7).Value = 17
$ws.Cells.Item(5,8).Value = 52788
$ws.Cells.Item(7,2).Value = 457656
$ws.Cells.Item(7,3).Value = 1541
$ws.Cells.Item(7,4).Value = 259143
$ws.Cells.Item(7,5).Value = 184008
$ws.Cells.Item(7,7).Value = 22
$ws.Cells.Item(7,8).Value = 14505
$ws.Cells.Item(13,2).Value = 212501
$ws.Cells.Item(13,3).Value = 2531
$ws.Cells.Item(13,4).Value = 172096
$ws.Cells.Item(13,5).Value = 30409
$ws.Cells.Item(13,7).Value = 133
$ws.Cells.Item(13,8).Value = 9996
$ws.Cells.Item(14,2).Value = 192827
$ws.Cells.Item(14,3).Value = 49
$ws.Cells.Item(14,5).Value = 8141
$ws.Cells.Item(23,2).Value = 90778
$ws.Cells.Item(23,3).Value = 1199
$ws.Cells.Item(23,4).Value = 73083
$ws.Cells.Item(23,5).Value = 17591
$ws.Cells.Item(23,7).Value = 5
$ws.Cells.Item(23,8).Value = 104
$ws.Cells.Item(28,2).Value = 59945
$ws.Cells.Item(28,3).Value = 458
$ws.Cells.Item(28,4).Value = 40136
$ws.Cells.Item(28,5).Value = 19447
$ws.Cells.Item(28,7).Value = 5
$ws.Cells.Item(28,8).Value = 362
$ws.Cells.Item(36,2).Value = 41879
$ws.Cells.Item(36,3).Value = 846
$ws.Cells.Item(36,4).Value = 32809
$ws.Cells.Item(36,5).Value = 8733
$ws.Cells.Item(36,7).Value = 3
$ws.Cells.Item(36,8).Value = 337
$ws.Cells.Item(42,2).Value = 32295
$ws.Cells.Item(42,3).Value = 470
$ws.Cells.Item(42,4).Value = 8656
$ws.Cells.Item(42,5).Value = 22435
$ws.Cells.Item(42,7).Value = 18
$ws.Cells.Item(42,8).Value = 1204
$ws.Cells.Item(43,2).Value = 31376
$ws.Cells.Item(43,3).Value = 44
$ws.Cells.Item(43,5).Value = 418
$ws.Cells.Item(43,7).Value = 2
$ws.Cells.Item(43,8).Value = 1958
$ws.Cells.Item(67,1).Value = "Nepal"
$ws.Cells.Item(67,2).Value = 10728
$ws.Cells.Item(67,3).Value = 629
$ws.Cells.Item(67,4).Value = 2338
$ws.Cells.Item(67,5).Value = 8366
$ws.Cells.Item(67,8).Value = 24
$ws.Cells.Item(68,1).Value = "Marruecos"
$ws.Cells.Item(68,2).Value = 10693
$ws.Cells.Item(68,3).Value = 349
$ws.Cells.Item(68,4).Value = 8426
$ws.Cells.Item(68,5).Value = 2053
$ws.Cells.Item(68,7).Value = 0
$ws.Cells.Item(68,8).Value = 214
$ws.Cells.Item(69,1).Value = "Chequia"
$ws.Cells.Item(69,2).Value = 10651
$ws.Cells.Item(69,3).Value = 1
$ws.Cells.Item(69,4).Value = 7559
$ws.Cells.Item(69,5).Value = 2752
$ws.Cells.Item(69,7).Value = 1
$ws.Cells.Item(69,8).Value = 340
$ws.Cells.Item(75,4).Value = 6600
$ws.Cells.Item(75,5).Value = 240
$ws.Cells.Item(76,4).Value = 4588
$ws.Cells.Item(76,5).Value = 2148
$ws.Cells.Item(107,2).Value = 1998
$ws.Cells.Item(107,3).Value = 7
$ws.Cells.Item(107,5).Value = 425
$ws.Cells.Item(113,2).Value = 1787
$ws.Cells.Item(113,3).Value = 63
$ws.Cells.Item(113,4).Value = 779
$ws.Cells.Item(113,5).Value = 992
$ws.Cells.Item(113,7).Value = 1
$ws.Cells.Item(113,8).Value = 16
$ws.Cells.Item(139,1).Value = "Uganda"
$ws.Cells.Item(139,2).Value = 805
$ws.Cells.Item(139,3).Value = 8
$ws.Cells.Item(139,4).Value = 717
$ws.Cells.Item(139,5).Value = 88
$ws.Cells.Item(139,8).Value = 0
$ws.Cells.Item(140,1).Value = "Malaui"
$ws.Cells.Item(140,2).Value = 803
$ws.Cells.Item(140,4).Value = 258
$ws.Cells.Item(140,5).Value = 534
$ws.Cells.Item(140,8).Value = 11
$ws.Cells.Item(141,1).Value = "Ruanda"
$ws.Cells.Item(141,2).Value = 798
$ws.Cells.Item(141,4).Value = 371
$ws.Cells.Item(141,5).Value = 425
$ws.Cells.Item(141,8).Value = 2
$ws.Cells.Item(145,4).Value = 647
$ws.Cells.Item(145,5).Value = 9
$ws.Cells.Item(148,4).Value = 624
$ws.Cells.Item(148,5).Value = 32
$ws.Cells.Item(157,2).Value = 352
$ws.Cells.Item(157,3).Value = 3
$ws.Cells.Item(157,5).Value = 23
$ws.Cells.Item(202,1).Value = "Dominica"
$ws.Cells.Item(203,1).Value = "Fiyi"
$ws.Cells.Item(211,1).Value = "Montserrat"
$ws.Cells.Item(211,4).Value = 10
$ws.Cells.Item(211,8).Value = 1
$ws.Cells.Item(212,1).Value = "Seychelles"
$ws.Cells.Item(212,4).Value = 11
$ws.Cells.Item(212,8).Value = 0
